$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2-15) held "D0(W) " (with a trailing space) - trim it to "D0(W)"
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "D0(W)"
}

# Move the active selection from G19 to D18
$ws.Range("D18").Select()
